$d = $word.ActiveDocument

# Update the Phase I SBIR/STTR grant amount from $256,000 to $275,000
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$replaced = $find.Execute("256,000", $true, $true, $false, $false, $false, $true, 1, $false, "275,000", 2)

if (-not $replaced) {
    throw "Could not find '256,000' to replace with '275,000'"
}

Write-Output "Replaced: $replaced"
